$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''304.04'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''0.15%'
$ws.Range('E2').Style = 'Normal'
$ws.Range('G2').Value = '''16'
$ws.Range('G2').Style = 'Normal'
$ws.Range('D3').Value = '''35.84'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''-3.64%'
$ws.Range('E3').Style = 'Normal'
$ws.Range('G3').Value = '''16'
$ws.Range('G3').Style = 'Normal'
$ws.Range('D4').Value = '''5.052'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''0.51%'
$ws.Range('E4').Style = 'Normal'
$ws.Range('G4').Value = '''16'
$ws.Range('G4').Style = 'Normal'
$ws.Range('D5').Value = '''0.07882'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''0.96%'
$ws.Range('E5').Style = 'Normal'
$ws.Range('G5').Value = '''16'
$ws.Range('G5').Style = 'Normal'
$ws.Range('D6').Value = '''2.129'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''-4.36%'
$ws.Range('E6').Style = 'Normal'
$ws.Range('G6').Value = '''16'
$ws.Range('G6').Style = 'Normal'
$ws.Range('B7').Value = 'KuCoinToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range('D7').Value = '''7.922'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''-0.89%'
$ws.Range('E7').Style = 'Normal'
$ws.Range('G7').Value = '''16'
$ws.Range('G7').Style = 'Normal'
$ws.Range('B8').Value = 'MXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D8').Value = '''0.9218'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''-0.46%'
$ws.Range('E8').Style = 'Normal'
$ws.Range('G8').Value = '''16'
$ws.Range('G8').Style = 'Normal'
$ws.Range('B9').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C9').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D9').Value = '''0.09764'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''-1.36%'
$ws.Range('E9').Style = 'Normal'
$ws.Range('G9').Value = '''16'
$ws.Range('G9').Style = 'Normal'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D10').Value = '''0.1855'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''-1.59%'
$ws.Range('E10').Style = 'Normal'
$ws.Range('G10').Value = '''16'
$ws.Range('G10').Style = 'Normal'
$ws.Range('B11').Value = 'MandalaExchangeToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D11').Value = '''0.08629'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''0.18%'
$ws.Range('E11').Style = 'Normal'
$ws.Range('G11').Value = '''16'
$ws.Range('G11').Style = 'Normal'
$ws.Range('B12').Value = 'BitrueCoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').Value = '''0.03601'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''-1.14%'
$ws.Range('E12').Style = 'Normal'
$ws.Range('G12').Value = '''16'
$ws.Range('G12').Style = 'Normal'
$ws.Range('B13').Value = 'BitMartToken'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').Value = '''0.09944'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''0.04%'
$ws.Range('E13').Style = 'Normal'
$ws.Range('G13').Value = '''16'
$ws.Range('G13').Style = 'Normal'
$ws.Range('B14').Value = 'BitForexToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').Value = '''0.001441'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''-2.06%'
$ws.Range('E14').Style = 'Normal'
$ws.Range('G14').Value = '''16'
$ws.Range('G14').Style = 'Normal'
$ws.Range('B15').Value = 'TigerCash'
$ws.Range('C15').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D15').Value = '''0.005716'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''0.76%'
$ws.Range('E15').Style = 'Normal'
$ws.Range('G15').Value = '''16'
$ws.Range('G15').Style = 'Normal'
$ws.Range('B16').Value = 'LEO'
$ws.Range('C16').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D16').Value = '''3.472'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''0.45%'
$ws.Range('E16').Style = 'Normal'
$ws.Range('G16').Value = '''16'
$ws.Range('G16').Style = 'Normal'
$ws.Range('B17').Value = 'GateToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D17').Value = '''4.127'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''2.46%'
$ws.Range('E17').Style = 'Normal'
$ws.Range('G17').Value = '''16'
$ws.Range('G17').Style = 'Normal'
$ws.Range('E18').Value = '''15.93%'
$ws.Range('E18').Style = 'Normal'
$ws.Range('G18').Value = '''16'
$ws.Range('G18').Style = 'Normal'
$ws.Range('E19').Value = '''-1.72%'
$ws.Range('E19').Style = 'Normal'
$ws.Range('G19').Value = '''16'
$ws.Range('G19').Style = 'Normal'
$ws.Range('D20').Value = '''0.1345'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''2.43%'
$ws.Range('E20').Style = 'Normal'
$ws.Range('G20').Value = '''16'
$ws.Range('G20').Style = 'Normal'
$ws.Range('D21').Value = '''5.152'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''8.02%'
$ws.Range('E21').Style = 'Normal'
$ws.Range('G21').Value = '''16'
$ws.Range('G21').Style = 'Normal'
$ws.Range('E22').Value = '''0.09%'
$ws.Range('E22').Style = 'Normal'
$ws.Range('G22').Value = '''16'
$ws.Range('G22').Style = 'Normal'
$ws.Range('D23').Value = '''0.04565'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''-0.78%'
$ws.Range('E23').Style = 'Normal'
$ws.Range('G23').Value = '''16'
$ws.Range('G23').Style = 'Normal'
$ws.Range('D24').Value = '''0.001231'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''-1.82%'
$ws.Range('E24').Style = 'Normal'
$ws.Range('G24').Value = '''16'
$ws.Range('G24').Style = 'Normal'
$ws.Range('D25').Value = '''0.004819'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '''-2.97%'
$ws.Range('E25').Style = 'Normal'
$ws.Range('G25').Value = '''16'
$ws.Range('G25').Style = 'Normal'
$ws.Range('D26').Value = '''0.0001301'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''-7.32%'
$ws.Range('E26').Style = 'Normal'
$ws.Range('G26').Value = '''16'
$ws.Range('G26').Style = 'Normal'
$ws.Range('E27').Value = '''74.42%'
$ws.Range('E27').Style = 'Normal'
$ws.Range('G27').Value = '''16'
$ws.Range('G27').Style = 'Normal'
$ws.Range('G28').Value = '''16'
$ws.Range('G28').Style = 'Normal'
$ws.Range('G29').Value = '''16'
$ws.Range('G29').Style = 'Normal'
$ws.Range('G30').Value = '''16'
$ws.Range('G30').Style = 'Normal'
$ws.Range('G31').Value = '''16'
$ws.Range('G31').Style = 'Normal'
$ws.Range('G32').Value = '''16'
$ws.Range('G32').Style = 'Normal'
$ws.Range('G33').Value = '''16'
$ws.Range('G33').Style = 'Normal'
$ws.Range('G34').Value = '''16'
$ws.Range('G34').Style = 'Normal'
$ws.Range('G35').Value = '''16'
$ws.Range('G35').Style = 'Normal'
$ws.Range('G36').Value = '''16'
$ws.Range('G36').Style = 'Normal'
$ws.Range('G37').Value = '''16'
$ws.Range('G37').Style = 'Normal'
$ws.Range('G38').Value = '''16'
$ws.Range('G38').Style = 'Normal'
$ws.Range('D39').Value = '''0.01878'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''2.36%'
$ws.Range('E39').Style = 'Normal'
$ws.Range('G39').Value = '''16'
$ws.Range('G39').Style = 'Normal'
$ws.Range('D40').Value = '''0.04707'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''-1.06%'
$ws.Range('E40').Style = 'Normal'
$ws.Range('G40').Value = '''16'
$ws.Range('G40').Style = 'Normal'
$ws.Range('D41').Value = '''0.007811'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''-2.13%'
$ws.Range('E41').Style = 'Normal'
$ws.Range('G41').Value = '''16'
$ws.Range('G41').Style = 'Normal'
$ws.Range('D42').Value = '''0.1389'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''-0.91%'
$ws.Range('E42').Style = 'Normal'
$ws.Range('G42').Value = '''16'
$ws.Range('G42').Style = 'Normal'
$ws.Range('D43').Value = '''0.007751'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''1.98%'
$ws.Range('E43').Style = 'Normal'
$ws.Range('G43').Value = '''16'
$ws.Range('G43').Style = 'Normal'
$ws.Range('D44').Value = '''0.002221'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''-1.08%'
$ws.Range('E44').Style = 'Normal'
$ws.Range('G44').Value = '''16'
$ws.Range('G44').Style = 'Normal'
$ws.Range('D45').Value = '''0.01137'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''8.95%'
$ws.Range('E45').Style = 'Normal'
$ws.Range('G45').Value = '''16'
$ws.Range('G45').Style = 'Normal'
$ws.Range('D46').Value = '''0.00006384'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''1.35%'
$ws.Range('E46').Style = 'Normal'
$ws.Range('G46').Value = '''16'
$ws.Range('G46').Style = 'Normal'
$ws.Range('E47').Value = '''-0.20%'
$ws.Range('E47').Style = 'Normal'
$ws.Range('G47').Value = '''16'
$ws.Range('G47').Style = 'Normal'
$ws.Range('E48').Value = '''-0.08%'
$ws.Range('E48').Style = 'Normal'
$ws.Range('G48').Value = '''16'
$ws.Range('G48').Style = 'Normal'
$ws.Range('D49').Value = '''51.86'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''68.48%'
$ws.Range('E49').Style = 'Normal'
$ws.Range('G49').Value = '''16'
$ws.Range('G49').Style = 'Normal'
$ws.Range('E50').Value = '''-29.49%'
$ws.Range('E50').Style = 'Normal'
$ws.Range('G50').Value = '''16'
$ws.Range('G50').Style = 'Normal'
$ws.Range('E51').Value = '''-0.20%'
$ws.Range('E51').Style = 'Normal'
$ws.Range('G51').Value = '''16'
$ws.Range('G51').Style = 'Normal'
